$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3915
$ws.Range("I2").Value = 1552.3334
$ws.Range("J2").Value = 4927.5713
$ws.Range("K2").Value = 1552.3334
$ws.Range("L2").Value = 4927.5713
$ws.Range("M2").Value = -1439.3334
$ws.Range("N2").Value = -5153.5713
$ws.Range("H42").Value = 3773.125
$ws.Range("J42").Value = 7512.5
$ws.Range("L42").Value = 22537.5
$ws.Range("N42").Value = -22997.5
$ws.Range("H104").Value = 212
$ws.Range("I104").Value = 212
$ws.Range("K104").Value = 636
$ws.Range("M104").Value = 1111
$ws.Range("H127").Value = 1674.25
$ws.Range("I127").Value = 1232.3334
$ws.Range("K127").Value = 3697.0002
$ws.Range("M127").Value = 1262.9998
$ws.Range("H132").Value = 2250.8696
$ws.Range("I132").Value = 2250.8696
$ws.Range("K132").Value = 6752.6088
$ws.Range("M132").Value = -4222.6088
$ws.Range("H137").Value = 2962.111
$ws.Range("J137").Value = 2986.8333
$ws.Range("L137").Value = 8960.499899999999
$ws.Range("N137").Value = -14060.4999
$ws.Range("H141").Value = 5054.8
$ws.Range("I141").Value = 5054.8
$ws.Range("K141").Value = 15164.4
$ws.Range("M141").Value = -9984.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3425.4285
$ws.Range("I32").Value = 3346.7
$ws.Range("K32").Value = 3346.7
$ws.Range("M32").Value = -3059.7
$ws.Range("H37").Value = 18122.5
$ws.Range("J37").Value = 19997.143
$ws.Range("L37").Value = 19997.143
$ws.Range("N37").Value = -20543.143
$ws.Range("H55").Value = 27855.715
$ws.Range("J55").Value = 27855.715
$ws.Range("L55").Value = 27855.715
$ws.Range("N55").Value = -28485.715
$ws.Range("H74").Value = 1461.1818
$ws.Range("I74").Value = 1461.1818
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1461.1818
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -587.1818000000001
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 1461.1818
$ws.Range("I77").Value = 1461.1818
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 7305.909000000001
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -2937.909000000001
$ws.Range("N77").ClearContents()
$ws.Range("H80").Value = 39998.332
$ws.Range("J80").Value = 39998.125
$ws.Range("L80").Value = 39998.125
$ws.Range("N80").Value = -41994.125
$ws.Range("H83").Value = 39998.332
$ws.Range("J83").Value = 39998.125
$ws.Range("L83").Value = 119994.375
$ws.Range("N83").Value = -129978.375
$ws.Range("H97").Value = 772.9091
$ws.Range("I97").Value = 690.2
$ws.Range("J97").Value = 1600
$ws.Range("K97").Value = 690.2
$ws.Range("L97").Value = 1600
$ws.Range("M97").Value = -194.2
$ws.Range("N97").Value = -2592
$ws.Range("H102").Value = 1327.2727
$ws.Range("I102").Value = 1177.7778
$ws.Range("K102").Value = 1177.7778
$ws.Range("M102").Value = 444.2221999999999
$ws.Range("H122").Value = 3076.8572
$ws.Range("I122").Value = 3076.8572
$ws.Range("K122").Value = 9230.571599999999
$ws.Range("M122").Value = -6780.571599999999
$ws.Range("H132").Value = 3500
$ws.Range("I132").Value = 3500
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10500
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7970
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5778.5713
$ws.Range("I86").Value = 5129.5713
$ws.Range("K86").Value = 5129.5713
$ws.Range("M86").Value = -4006.5713
$ws.Range("H89").Value = 5778.5713
$ws.Range("I89").Value = 5129.5713
$ws.Range("K89").Value = 25647.8565
$ws.Range("M89").Value = -20031.8565
$ws.Range("H94").Value = 1529.8636
$ws.Range("I94").Value = 1428.2632
$ws.Range("J94").Value = 2173.3333
$ws.Range("K94").Value = 1428.2632
$ws.Range("L94").Value = 2173.3333
$ws.Range("M94").Value = -977.2632000000001
$ws.Range("N94").Value = -3075.3333
$ws.Range("H105").Value = 2675.25
$ws.Range("I105").Value = 2542.1667
$ws.Range("J105").Value = 3074.5
$ws.Range("K105").Value = 2542.1667
$ws.Range("L105").Value = 3074.5
$ws.Range("M105").Value = -795.1667000000002
$ws.Range("N105").Value = -6568.5
$ws.Range("H107").Value = 1060.5625
$ws.Range("I107").Value = 844.1539
$ws.Range("J107").Value = 1998.3334
$ws.Range("K107").Value = 844.1539
$ws.Range("L107").Value = 1998.3334
$ws.Range("M107").Value = 1075.8461
$ws.Range("N107").Value = -5838.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 4665
$ws.Range("J15").Value = 4665
$ws.Range("L15").Value = 4665
$ws.Range("N15").Value = -5005
$ws.Range("H31").Value = 2500.7144
$ws.Range("I31").Value = 1938.3334
$ws.Range("K31").Value = 1938.3334
$ws.Range("M31").Value = -1643.3334
$ws.Range("H34").Value = 2500.7144
$ws.Range("I34").Value = 1938.3334
$ws.Range("K34").Value = 1938.3334
$ws.Range("M34").Value = -1736.3334
$ws.Range("H58").Value = 4671.8335
$ws.Range("I58").Value = 4053.7693
$ws.Range("K58").Value = 4053.7693
$ws.Range("M58").Value = -3850.7693
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H74").Value = 39998.57
$ws.Range("J74").Value = 39998.57
$ws.Range("L74").Value = 39998.57
$ws.Range("N74").Value = -41746.57
$ws.Range("H77").Value = 39998.57
$ws.Range("J77").Value = 39998.57
$ws.Range("L77").Value = 119995.71
$ws.Range("N77").Value = -128731.71
$ws.Range("H105").Value = 1713.1666
$ws.Range("I105").Value = 1426.6666
$ws.Range("K105").Value = 1426.6666
$ws.Range("M105").Value = 320.3334
$ws.Range("H132").Value = 2063.5715
$ws.Range("I132").Value = 1979
$ws.Range("J132").Value = 2275
$ws.Range("K132").Value = 5937
$ws.Range("L132").Value = 6825
$ws.Range("M132").Value = -3407
$ws.Range("N132").Value = -11885
$ws.Range("H134").Value = 2168.5908
$ws.Range("I134").Value = 1782.2354
$ws.Range("K134").Value = 5346.706200000001
$ws.Range("M134").Value = -2811.706200000001
$ws.Range("H136").Value = 4671.8335
$ws.Range("I136").Value = 4053.7693
$ws.Range("K136").Value = 12161.3079
$ws.Range("M136").Value = -9611.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H107").Value = 797.5
$ws.Range("I107").Value = 1030
$ws.Range("K107").Value = 1030
$ws.Range("M107").Value = 890
$ws.Range("H126").Value = 16661.375
$ws.Range("I126").Value = 14755.857
$ws.Range("K126").Value = 44267.571
$ws.Range("M126").Value = -41797.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 935.1111
$ws.Range("I22").Value = 825.63635
$ws.Range("J22").Value = 1107.1428
$ws.Range("K22").Value = 825.63635
$ws.Range("L22").Value = 1107.1428
$ws.Range("M22").Value = -530.63635
$ws.Range("N22").Value = -1697.1428
$ws.Range("H27").Value = 935.1111
$ws.Range("I27").Value = 825.63635
$ws.Range("J27").Value = 1107.1428
$ws.Range("K27").Value = 825.63635
$ws.Range("L27").Value = 1107.1428
$ws.Range("M27").Value = -718.63635
$ws.Range("N27").Value = -1321.1428
$ws.Range("H61").Value = 2331.6667
$ws.Range("I61").Value = 2331.6667
$ws.Range("K61").Value = 2331.6667
$ws.Range("M61").Value = -2129.6667
$ws.Range("H68").Value = 30214
$ws.Range("I68").Value = 1699.8
$ws.Range("J68").Value = 101499.5
$ws.Range("K68").Value = 1699.8
$ws.Range("L68").Value = 101499.5
$ws.Range("M68").Value = -950.8
$ws.Range("N68").Value = -102997.5
$ws.Range("H71").Value = 30214
$ws.Range("I71").Value = 1699.8
$ws.Range("J71").Value = 101499.5
$ws.Range("K71").Value = 8499
$ws.Range("L71").Value = 507497.5
$ws.Range("M71").Value = -4755
$ws.Range("N71").Value = -514985.5
$ws.Range("H113").Value = 2331.6667
$ws.Range("I113").Value = 2331.6667
$ws.Range("K113").Value = 2331.6667
$ws.Range("M113").Value = -161.6667000000002
$ws.Range("H122").Value = 3456.7273
$ws.Range("I122").Value = 3224.3333
$ws.Range("K122").Value = 9672.999899999999
$ws.Range("M122").Value = -7222.999899999999
$ws.Range("H132").Value = 3197.8
$ws.Range("I132").Value = 1095.1111
$ws.Range("J132").Value = 4918.1816
$ws.Range("K132").Value = 3285.3333
$ws.Range("L132").Value = 14754.5448
$ws.Range("M132").Value = -755.3333000000002
$ws.Range("N132").Value = -19814.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 29975
$ws.Range("J54").Value = 29975
$ws.Range("L54").Value = 29975
$ws.Range("N54").Value = -31015
$ws.Range("H113").Value = 602.55554
$ws.Range("I113").Value = 466
$ws.Range("K113").Value = 1398
$ws.Range("M113").Value = 772
$ws.Range("H132").Value = 1115.8636
$ws.Range("I132").Value = 1103.6111
$ws.Range("K132").Value = 3310.8333
$ws.Range("M132").Value = -780.8333000000002
